$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.26415479183197
$ws.Range("B1").Value = 2.770501375198364
$ws.Range("C1").Value = 8.83018970489502
$ws.Range("D1").Value = 2.036558628082275
$ws.Range("E1").Value = 1.130433559417725
